# Update Week (row 2) target depth data on both the OFF and DEF sheets
# to reflect newly logged Week 16 stats (and resulting season-sim recalculation).

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 191
$wsOff.Range("C2").Value = 135
$wsOff.Range("D2").Value = 50
$wsOff.Range("E2").Value = 29

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 185
$wsDef.Range("C2").Value = 116
$wsDef.Range("D2").Value = 48
$wsDef.Range("E2").Value = 29
$wsDef.Range("F2").Value = 6
$wsDef.Range("G2").Value = 3
